# Recompute the "Total" column (AB) so it no longer includes the
# erroneous 4x multiplier: new value = old value / 4
# (equivalently, Total = sum of the per-asset columns B:AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 28)  # Column AB = 28
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val / 4
    }
}
